$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 208: 28-10-2021
$ws.Range("A208").Value = "28-10-2021"
$ws.Range("B208").Value = 3376000
$ws.Range("C208").Value = 87841

# Row 209: 29-10-2021
$ws.Range("A209").Value = "29-10-2021"
$ws.Range("B209").Value = 3893000
$ws.Range("C209").Value = 24000

# Row 210: 02-11-2021
# Use a leading apostrophe to force Excel to store the value as literal
# text instead of auto-converting it to a date serial number, then clear
# the resulting cell formatting so no extra number format is applied.
$ws.Range("A210").Value = "'02-11-2021"
$ws.Range("A210").ClearFormats()
$ws.Range("B210").Value = 4433800
$ws.Range("C210").Value = 0
